# feat: add 2022-Q1 data
#
# Target layout after edit:
#   2020-Q4 (unchanged)
#   2021-Q4 (unchanged)
#   2022-Q1 (new fund-holding sheet, inserted before "总计")
#   总计    (existing summary sheet, gets a new first data row for 2022-Q1)
#
# NOTE: worksheet object references returned by Worksheets.Item(...) seem to
# track POSITION rather than stable identity in this host, so after any
# operation that changes sheet order/count (Copy/Add/Delete/Move) we
# re-fetch the worksheets we still need by name/index instead of reusing
# older variables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the existing "总计" sheet; the copy lands right after the
#    original. The ORIGINAL (3rd sheet) becomes the new "2022-Q1"
#    fund-holdings sheet, and the COPY (4th sheet) becomes the refreshed
#    "总计" summary sheet that now includes the 2022-Q1 row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Copy($null, $totalSheet)

$q1Sheet = $wb.Worksheets.Item(3)
$q1Sheet.Name = "2022-Q1"

$newTotalSheet = $wb.Worksheets.Item(4)
$newTotalSheet.Name = "总计"

# ---------------------------------------------------------------------
# 2) Rebuild "2022-Q1": it is currently still shaped like the old "总计"
#    sheet (4 columns, date/count/value). Turn it into the 8-column fund
#    holdings layout used by "2020-Q4" / "2021-Q4".
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# extend the bold/centered/bordered header style (already on B1:D1) across
# the new E1:H1 header cells
$q1Sheet.Range("D1").Copy()
$q1Sheet.Range("E1:H1").PasteSpecial(-4122)

$q1Sheet.Cells.Item(1, 2).Value = "基金代码"
$q1Sheet.Cells.Item(1, 3).Value = "基金名称"
$q1Sheet.Cells.Item(1, 4).Value = "基金规模"
$q1Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q1Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q1Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1Sheet.Cells.Item(1, 8).Value = "仓位排名"

# force text formatting on the data columns so fund codes keep their
# leading zeros and the numeric-looking values stay text (matches the
# rest of the workbook)
$q1Sheet.Range("B2:G3").NumberFormat = "@"

$q1Sheet.Cells.Item(2, 1).Value = 0
$q1Sheet.Cells.Item(2, 2).Value = "009686"
$q1Sheet.Cells.Item(2, 3).Value = "华夏磐利一年定期开放混合A"
$q1Sheet.Cells.Item(2, 4).Value = "16.02"
$q1Sheet.Cells.Item(2, 5).Value = "93.69"
$q1Sheet.Cells.Item(2, 6).Value = "4.01"
$q1Sheet.Cells.Item(2, 7).Value = "0.6424"
$q1Sheet.Cells.Item(2, 8).Value = 7

$q1Sheet.Cells.Item(3, 1).Value = 1
$q1Sheet.Cells.Item(3, 2).Value = "009687"
$q1Sheet.Cells.Item(3, 3).Value = "华夏磐利一年定期开放混合C"
$q1Sheet.Cells.Item(3, 4).Value = "0.51"
$q1Sheet.Cells.Item(3, 5).Value = "93.69"
$q1Sheet.Cells.Item(3, 6).Value = "4.01"
$q1Sheet.Cells.Item(3, 7).Value = "0.0205"
$q1Sheet.Cells.Item(3, 8).Value = 7

# ---------------------------------------------------------------------
# 3) "总计" summary sheet: insert a new row right under the header and
#    fill it in with the 2022-Q1 totals. The two previously-existing
#    rows (2021-Q4, 2020-Q4) shift down.
# ---------------------------------------------------------------------
$newTotalSheet = $wb.Worksheets.Item("总计")

$newTotalSheet.Rows.Item(2).Insert()
$newTotalSheet.Rows.Item(2).ClearFormats()

$newTotalSheet.Cells.Item(2, 1).Value = 0
$newTotalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$newTotalSheet.Cells.Item(2, 3).Value = 2
$newTotalSheet.Cells.Item(2, 4).Value = 0.66

# renumber the index column for the rows that shifted down
$newTotalSheet.Cells.Item(3, 1).Value = 1
$newTotalSheet.Cells.Item(4, 1).Value = 2

# restore the index-column styling (border/bold/center) on the new row
$newTotalSheet.Range("A3").Copy()
$newTotalSheet.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4) Restore the original active sheet/selection (the diff doesn't touch
#    the workbook-level active tab, which pointed at "2020-Q4").
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$firstSheet.Activate()
$firstSheet.Range("A1").Select()
